$d = $word.ActiveDocument

# The "Write Up" title paragraph is the first paragraph in the document.
$titlePara = $d.Paragraphs(1)
$titleRange = $titlePara.Range

# Insert three new, empty paragraphs right after the title - these will
# become paragraphs 2, 3 and 4 - and then fill them in with the new
# write-up / SEO copy and the new "The JavaScript File" heading.
$titleRange.InsertParagraphAfter()
$titleRange.InsertParagraphAfter()
$titleRange.InsertParagraphAfter()

$introPara = $d.Paragraphs(2)
$introPara.Style = "Normal"
$introPara.Range.Text = "This week, we will be adding some JavaScript to our project. If we want our hamburger menu to be functional and actually do something kind-of cool, when we click on it, we will need to add some JavaScript."

$teaserPara = $d.Paragraphs(3)
$teaserPara.Style = "Normal"
$teaserPara.Range.Text = "So, if this tends to pique your interest, then please join us for our brand-new article this week entitled:"

$headingPara = $d.Paragraphs(4)
$headingPara.Style = "Heading1"
$headingPara.Range.Text = "The JavaScript File"
